$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Initial Position Single Thread")

# --- New shared string used by the commit note column (P) ---
$noteText = "avoided obvious ray attacks calculation on move validation"

# --- Normalise the style of the existing "note" cells so they match the
#     default style used by the rest of that column (P76 already uses it).
#     A plain value re-write resets the style to the default (0); applying
#     the "Normal" cell style nudges it back to the alignment/protection
#     flavoured default that the rest of the sheet uses. ---
$ws.Range("P77").Value = "added functions map for destination calculation"
$ws.Range("P77").Style = "Normal"

$ws.Range("P80").Value = "added pawn attacks functional maps"
$ws.Range("P80").Style = "Normal"

$ws.Range("P81").Value = "replaced unsigned char with unsigned int"
$ws.Range("P81").Style = "Normal"

# --- New benchmark block (rows 84-86), mirrors the row 80-82 block with
#     new run data but the same formula pattern shifted down by 4 rows. ---

# Row 84
$ws.Range("A84").Value = 45970
$ws.Range("C84").Value = 4
$ws.Range("D84").Value = 206603
$ws.Range("E84").Value = 206
$ws.Range("F84").Formula = "=D84/E84*1000"
$ws.Range("G84").Formula = "=(E80-E84)/E80"
$ws.Range("H84").Formula = "=(F84-80000000)/80000000"
$ws.Range("I84").Value = 4
$ws.Range("J84").Value = 197281
$ws.Range("K84").Value = 6
$ws.Range("L84").Formula = "=J84/K84*1000"
$ws.Range("M84").Formula = "=(K80-K84)/K80"
$ws.Range("N84").Formula = "=(L84-80000000)/80000000"
$ws.Range("P84").Value = $noteText

# Row 85
$ws.Range("C85").Value = 5
$ws.Range("D85").Value = 5072212
$ws.Range("E85").Value = 5048
$ws.Range("F85").Formula = "=D85/E85*1000"
$ws.Range("G85").Formula = "=(E81-E85)/E81"
$ws.Range("H85").Formula = "=(F85-80000000)/80000000"
$ws.Range("I85").Value = 5
$ws.Range("J85").Value = 4880523
$ws.Range("K85").Value = 155
$ws.Range("L85").Formula = "=J85/K85*1000"
$ws.Range("M85").Formula = "=(K81-K85)/K81"
$ws.Range("N85").Formula = "=(L85-80000000)/80000000"

# Row 86
$ws.Range("I86").Value = 6
$ws.Range("J86").Value = 119060324
$ws.Range("K86").Value = 3790
$ws.Range("L86").Formula = "=J86/K86*1000"
$ws.Range("M86").Formula = "=(K82-K86)/K82"
$ws.Range("N86").Formula = "=(L86-80000000)/80000000"

# --- Formatting for the new block: reuse the same direct formatting as
#     the block above it (rows 80-82), column by column. ---
foreach ($col in @("A","C","D","E","F","G","H","I","J","K","L","M","N")) {
    $ws.Range("$col 84".Replace(" ", "")).NumberFormat = $ws.Range("$col 80".Replace(" ", "")).NumberFormat
    $ws.Range("$col 85".Replace(" ", "")).NumberFormat = $ws.Range("$col 81".Replace(" ", "")).NumberFormat
    $ws.Range("$col 86".Replace(" ", "")).NumberFormat = $ws.Range("$col 82".Replace(" ", "")).NumberFormat
}
$ws.Range("P84").Style = "Normal"

# --- Book-keeping: dimension grows, selection lands on the last edited note cell ---
$ws.Range("P85").Select()

$wb.Save()
